# Add the new "medbert1" results sheet as the last tab, after "basebert1",
# matching the new <sheet> entry appended to workbook.xml in the diff.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "medbert1"

# --- Header row: B1:N1 text labels ---
$ws.Range("B1").Value = "Fold"
$ws.Range("C1").Value = "Version"
$ws.Range("D1").Value = "Epoch"
$ws.Range("E1").Value = "Recall"
$ws.Range("F1").Value = "Precision"
$ws.Range("G1").Value = "Accuracy"
$ws.Range("H1").Value = "Fbeta"
$ws.Range("I1").Value = "Best Recall"
$ws.Range("J1").Value = "Best Precision"
$ws.Range("K1").Value = "Best Threshold"
$ws.Range("L1").Value = "False Neg(0.5)"
$ws.Range("M1").Value = "False Pos(0.5)"
$ws.Range("N1").Value = "Val loss"

# --- Column A (fold index), numeric ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# --- Plain text cells (fold name, version timestamp) ---
$ws.Range("B2").Value = "fold_0"
$ws.Range("C2").Value = "18.09_11.56"
$ws.Range("B3").Value = "fold_1"
$ws.Range("C3").Value = "18.09_12.10"
$ws.Range("B4").Value = "fold_2"
$ws.Range("C4").Value = "18.09_12.24"
$ws.Range("B5").Value = "fold_3"
$ws.Range("C5").Value = "18.09_12.38"

# --- Numeric-looking metric cells stored as TEXT in the source (Recall/Precision/Accuracy/Fbeta) ---
# Force text number format per-cell first so the numeric-looking string is not auto-coerced to a number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("E2").Value = "0.64935064"
$ws.Range("F2").Value = "0.28089887"
$ws.Range("G2").Value = "0.8647469"
$ws.Range("H2").Value = "0.5144033"
$ws.Range("E3").Value = "0.8051948"
$ws.Range("F3").Value = "0.31958762"
$ws.Range("G3").Value = "0.87172776"
$ws.Range("H3").Value = "0.61752987"
$ws.Range("E4").Value = "0.8831169"
$ws.Range("F4").Value = "0.3090909"
$ws.Range("G4").Value = "0.8595113"
$ws.Range("H4").Value = "0.6439394"
$ws.Range("E5").Value = "0.8717949"
$ws.Range("F5").Value = "0.272"
$ws.Range("G5").Value = "0.83246076"
$ws.Range("H5").Value = "0.6049822"

# --- Numeric cells (epoch, thresholds, counts, val loss) ---
$ws.Range("D2").Value = 8
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0.5
$ws.Range("L2").Value = 27
$ws.Range("M2").Value = 128
$ws.Range("N2").Value = 0.6564610414206982
$ws.Range("D3").Value = 7
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.5
$ws.Range("L3").Value = 15
$ws.Range("M3").Value = 132
$ws.Range("N3").Value = 0.5011956257124742
$ws.Range("D4").Value = 7
$ws.Range("I4").Value = 0.961
$ws.Range("J4").Value = 0.2298
$ws.Range("K4").Value = 0.2207
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = 152
$ws.Range("N4").Value = 0.4681512216726939
$ws.Range("D5").Value = 7
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0.5
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 182
$ws.Range("N5").Value = 0.5822798783580462

# --- Formatting: bold + thin border + center/top align, matching basebert1's header/index style ---
$headerRange = $ws.Range("B1:N1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$indexRange = $ws.Range("A2:A5")
$indexRange.Font.Bold = $true
$indexRange.Borders.LineStyle = 1
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
